$wb = $excel.ActiveWorkbook
$wsConfig = $wb.Worksheets.Item("配置")
$wsSemantic = $wb.Worksheets.Item("语义分析")

# Update Ollama Host value: 192.168.96.246 -> 127.0.0.1
$wsConfig.Range("B3").Value = "127.0.0.1"

# Update Ollama Model value: wangshenzhi/llama3-8b-chinese-chat-ollama-q8 -> deepseek-r1:14b
$wsConfig.Range("B5").Value = "deepseek-r1:14b"

# Insert a new row 7 for the new option "显示语义分析打分及内容" (shows semantic
# analysis score & content), pushing the existing row 7 ("生成AI处理前后对比文档")
# down to row 8.
$wsConfig.Rows.Item(7).Insert()
$wsConfig.Range("A7").Value = "显示语义分析打分及内容"
$wsConfig.Range("B7").Value = "Y"

# Switch the active/selected sheet & selection to match the saved view state:
# "语义分析" is selected first (so it ends up with the non-active selection),
# then "配置" is activated last so it becomes the active tab.
$wsSemantic.Activate()
$wsSemantic.Range("E9").Select()

$wsConfig.Activate()
$wsConfig.Range("B3").Select()
